# Auto-generated edit script: updates cached market-price-derived leve profit
# figures (currentAveragePrice*, Leve*Price*, LeveProfit*) across the ALC, ARM,
# BSM, CRP, CUL, GSM, LTW, and WVR sheets, mirroring a scheduled price refresh.
$wb = $excel.ActiveWorkbook
$totalCells = 0

$ws = $wb.Worksheets.Item("ALC")
# Row 64
$ws.Range("H64").Value = 3831.111
$totalCells++

# Row 67
$ws.Range("H67").Value = 3831.111
$totalCells++

# Row 123
$ws.Range("H123").Value = 40945
$totalCells++
$ws.Range("J123").Value = 41890
$totalCells++
$ws.Range("L123").Value = 41890
$totalCells++
$ws.Range("N123").Value = -51690
$totalCells++

# Row 132
$ws.Range("H132").Value = 27031738
$totalCells++
$ws.Range("I132").Value = 34487836
$totalCells++
$ws.Range("J132").Value = 3376
$totalCells++
$ws.Range("K132").Value = 103463508
$totalCells++
$ws.Range("L132").Value = 10128
$totalCells++
$ws.Range("M132").Value = -103460978
$totalCells++
$ws.Range("N132").Value = -15188
$totalCells++

# Row 138
$ws.Range("H138").Value = 2578.88
$totalCells++
$ws.Range("I138").Value = 733
$totalCells++
$ws.Range("J138").Value = 2930.476
$totalCells++
$ws.Range("K138").Value = 2199
$totalCells++
$ws.Range("L138").Value = 8791.428
$totalCells++
$ws.Range("M138").Value = 2941
$totalCells++
$ws.Range("N138").Value = -19071.428
$totalCells++

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 3232.5
$totalCells++
$ws.Range("I45").Value = 3232.5
$totalCells++
$ws.Range("J45").Value = 0
$totalCells++
$ws.Range("K45").Value = 3232.5
$totalCells++
$ws.Range("L45").Value = 0
$totalCells++
$ws.Range("M45").Value = -2855.5
$totalCells++
$ws.Range("N45").ClearContents()
$totalCells++

# Row 61
$ws.Range("H61").Value = 1622.9412
$totalCells++
$ws.Range("I61").Value = 1545
$totalCells++
$ws.Range("J61").Value = 1986.6666
$totalCells++
$ws.Range("K61").Value = 1545
$totalCells++
$ws.Range("L61").Value = 1986.6666
$totalCells++
$ws.Range("M61").Value = -1333
$totalCells++
$ws.Range("N61").Value = -2410.6666
$totalCells++

# Row 97
$ws.Range("H97").Value = 1654.4615
$totalCells++
$ws.Range("I97").Value = 863.36365
$totalCells++
$ws.Range("K97").Value = 863.36365
$totalCells++
$ws.Range("M97").Value = -367.36365
$totalCells++

# Row 121
$ws.Range("H121").Value = 28849
$totalCells++
$ws.Range("J121").Value = 28849
$totalCells++
$ws.Range("L121").Value = 28849
$totalCells++
$ws.Range("N121").Value = -32343
$totalCells++

# Row 136
$ws.Range("H136").Value = 1622.9412
$totalCells++
$ws.Range("I136").Value = 1545
$totalCells++
$ws.Range("J136").Value = 1986.6666
$totalCells++
$ws.Range("K136").Value = 4635
$totalCells++
$ws.Range("L136").Value = 5959.9998
$totalCells++
$ws.Range("M136").Value = -2085
$totalCells++
$ws.Range("N136").Value = -11059.9998
$totalCells++

$ws = $wb.Worksheets.Item("BSM")
# Row 80
$ws.Range("H80").Value = 289.93332
$totalCells++
$ws.Range("I80").Value = 280
$totalCells++
$ws.Range("J80").Value = 294.9
$totalCells++
$ws.Range("K80").Value = 280
$totalCells++
$ws.Range("L80").Value = 294.9
$totalCells++
$ws.Range("M80").Value = 718
$totalCells++
$ws.Range("N80").Value = -2290.9
$totalCells++

# Row 83
$ws.Range("H83").Value = 289.93332
$totalCells++
$ws.Range("I83").Value = 280
$totalCells++
$ws.Range("J83").Value = 294.9
$totalCells++
$ws.Range("K83").Value = 1400
$totalCells++
$ws.Range("L83").Value = 1474.5
$totalCells++
$ws.Range("M83").Value = 3592
$totalCells++
$ws.Range("N83").Value = -11458.5
$totalCells++

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 5238.1514
$totalCells++
$ws.Range("I31").Value = 1128.15
$totalCells++
$ws.Range("J31").Value = 11561.23
$totalCells++
$ws.Range("K31").Value = 1128.15
$totalCells++
$ws.Range("L31").Value = 11561.23
$totalCells++
$ws.Range("M31").Value = -833.1500000000001
$totalCells++
$ws.Range("N31").Value = -12151.23
$totalCells++

# Row 34
$ws.Range("H34").Value = 5238.1514
$totalCells++
$ws.Range("I34").Value = 1128.15
$totalCells++
$ws.Range("J34").Value = 11561.23
$totalCells++
$ws.Range("K34").Value = 1128.15
$totalCells++
$ws.Range("L34").Value = 11561.23
$totalCells++
$ws.Range("M34").Value = -926.1500000000001
$totalCells++
$ws.Range("N34").Value = -11965.23
$totalCells++

# Row 38
$ws.Range("H38").Value = 19557.611
$totalCells++
$ws.Range("I38").Value = 1038
$totalCells++
$ws.Range("J38").Value = 20647
$totalCells++
$ws.Range("K38").Value = 1038
$totalCells++
$ws.Range("L38").Value = 20647
$totalCells++
$ws.Range("M38").Value = -661
$totalCells++
$ws.Range("N38").Value = -21401
$totalCells++

# Row 46
$ws.Range("H46").Value = 19557.611
$totalCells++
$ws.Range("I46").Value = 1038
$totalCells++
$ws.Range("J46").Value = 20647
$totalCells++
$ws.Range("K46").Value = 1038
$totalCells++
$ws.Range("L46").Value = 20647
$totalCells++
$ws.Range("M46").Value = -827
$totalCells++
$ws.Range("N46").Value = -21069
$totalCells++

# Row 50
$ws.Range("H50").Value = 29658.182
$totalCells++
$ws.Range("J50").Value = 29658.182
$totalCells++
$ws.Range("L50").Value = 29658.182
$totalCells++
$ws.Range("N50").Value = -30908.182
$totalCells++

# Row 57
$ws.Range("H57").Value = 45009.668
$totalCells++
$ws.Range("J57").Value = 45009.668
$totalCells++
$ws.Range("L57").Value = 45009.668
$totalCells++
$ws.Range("N57").Value = -46129.668
$totalCells++

# Row 68
$ws.Range("H68").Value = 57286.43
$totalCells++
$ws.Range("J68").Value = 57286.43
$totalCells++
$ws.Range("L68").Value = 57286.43
$totalCells++
$ws.Range("N68").Value = -58784.43
$totalCells++

# Row 71
$ws.Range("H71").Value = 57286.43
$totalCells++
$ws.Range("J71").Value = 57286.43
$totalCells++
$ws.Range("L71").Value = 171859.29
$totalCells++
$ws.Range("N71").Value = -179347.29
$totalCells++

# Row 86
$ws.Range("H86").Value = 5249.75
$totalCells++
$ws.Range("I86").Value = 7000
$totalCells++
$ws.Range("K86").Value = 7000
$totalCells++
$ws.Range("M86").Value = -5877
$totalCells++

# Row 89
$ws.Range("H89").Value = 5249.75
$totalCells++
$ws.Range("I89").Value = 7000
$totalCells++
$ws.Range("K89").Value = 35000
$totalCells++
$ws.Range("M89").Value = -29384
$totalCells++

# Row 123
$ws.Range("H123").Value = 40998.89
$totalCells++
$ws.Range("J123").Value = 40998.89
$totalCells++
$ws.Range("L123").Value = 40998.89
$totalCells++
$ws.Range("N123").Value = -50798.89
$totalCells++

# Row 132
$ws.Range("H132").Value = 2874.7896
$totalCells++
$ws.Range("I132").Value = 1774.8667
$totalCells++
$ws.Range("J132").Value = 6999.5
$totalCells++
$ws.Range("K132").Value = 5324.6001
$totalCells++
$ws.Range("L132").Value = 20998.5
$totalCells++
$ws.Range("M132").Value = -2794.6001
$totalCells++
$ws.Range("N132").Value = -26058.5
$totalCells++

$ws = $wb.Worksheets.Item("CUL")
# Row 34
$ws.Range("H34").Value = 10131.923
$totalCells++
$ws.Range("J34").Value = 7886.5586
$totalCells++
$ws.Range("L34").Value = 23659.6758
$totalCells++
$ws.Range("N34").Value = -23827.6758
$totalCells++

# Row 39
$ws.Range("H39").Value = 15776.947
$totalCells++
$ws.Range("J39").Value = 15776.947
$totalCells++
$ws.Range("L39").Value = 47330.841
$totalCells++
$ws.Range("N39").Value = -47918.841
$totalCells++

# Row 113
$ws.Range("H113").Value = 3572144
$totalCells++
$ws.Range("I113").Value = 588.7143
$totalCells++
$ws.Range("K113").Value = 1766.1429
$totalCells++
$ws.Range("M113").Value = 403.8571000000002
$totalCells++

$ws = $wb.Worksheets.Item("GSM")
# Row 123
$ws.Range("H123").Value = 10324.267
$totalCells++
$ws.Range("J123").Value = 10324.267
$totalCells++
$ws.Range("L123").Value = 10324.267
$totalCells++
$ws.Range("N123").Value = -15224.267
$totalCells++

# Row 126
$ws.Range("H126").Value = 3140.1
$totalCells++
$ws.Range("I126").Value = 2767.7
$totalCells++
$ws.Range("J126").Value = 4629.7
$totalCells++
$ws.Range("K126").Value = 8303.099999999999
$totalCells++
$ws.Range("L126").Value = 13889.1
$totalCells++
$ws.Range("M126").Value = -5833.099999999999
$totalCells++
$ws.Range("N126").Value = -18829.1
$totalCells++

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3974
$totalCells++
$ws.Range("I7").Value = 3394.45
$totalCells++
$ws.Range("J7").Value = 5027.727
$totalCells++
$ws.Range("K7").Value = 3394.45
$totalCells++
$ws.Range("L7").Value = 5027.727
$totalCells++
$ws.Range("M7").Value = -3282.45
$totalCells++
$ws.Range("N7").Value = -5251.727
$totalCells++

# Row 18
$ws.Range("H18").Value = 0
$totalCells++
$ws.Range("J18").Value = 0
$totalCells++
$ws.Range("L18").Value = 0
$totalCells++
$ws.Range("N18").ClearContents()
$totalCells++

# Row 20
$ws.Range("H20").Value = 7648.6875
$totalCells++
$ws.Range("J20").Value = 9770.817999999999
$totalCells++
$ws.Range("L20").Value = 9770.817999999999
$totalCells++
$ws.Range("N20").Value = -10222.818
$totalCells++

# Row 93
$ws.Range("H93").Value = 6175060.5
$totalCells++
$ws.Range("I93").Value = 8548737
$totalCells++
$ws.Range("J93").Value = 3500.8
$totalCells++
$ws.Range("K93").Value = 8548737
$totalCells++
$ws.Range("L93").Value = 3500.8
$totalCells++
$ws.Range("M93").Value = -8547489
$totalCells++
$ws.Range("N93").Value = -5996.8
$totalCells++

# Row 122
$ws.Range("H122").Value = 5747
$totalCells++
$ws.Range("I122").Value = 3661.6667
$totalCells++
$ws.Range("J122").Value = 8875
$totalCells++
$ws.Range("K122").Value = 10985.0001
$totalCells++
$ws.Range("L122").Value = 26625
$totalCells++
$ws.Range("M122").Value = -8535.000100000001
$totalCells++
$ws.Range("N122").Value = -31525
$totalCells++

# Row 126
$ws.Range("H126").Value = 3974
$totalCells++
$ws.Range("I126").Value = 3394.45
$totalCells++
$ws.Range("J126").Value = 5027.727
$totalCells++
$ws.Range("K126").Value = 10183.35
$totalCells++
$ws.Range("L126").Value = 15083.181
$totalCells++
$ws.Range("M126").Value = -7713.349999999999
$totalCells++
$ws.Range("N126").Value = -20023.181
$totalCells++

$ws = $wb.Worksheets.Item("WVR")
# Row 46
$ws.Range("H46").Value = 61626
$totalCells++
$ws.Range("J46").Value = 61626
$totalCells++
$ws.Range("L46").Value = 61626
$totalCells++
$ws.Range("N46").Value = -62088
$totalCells++

# Row 80
$ws.Range("H80").Value = 25200
$totalCells++
$ws.Range("J80").Value = 25200
$totalCells++
$ws.Range("L80").Value = 25200
$totalCells++
$ws.Range("N80").Value = -27196
$totalCells++

# Row 83
$ws.Range("H83").Value = 25200
$totalCells++
$ws.Range("J83").Value = 25200
$totalCells++
$ws.Range("L83").Value = 75600
$totalCells++
$ws.Range("N83").Value = -85584
$totalCells++

# Row 122
$ws.Range("H122").Value = 9375.9
$totalCells++
$ws.Range("I122").Value = 7400.5713
$totalCells++
$ws.Range("J122").Value = 13985
$totalCells++
$ws.Range("K122").Value = 22201.7139
$totalCells++
$ws.Range("L122").Value = 41955
$totalCells++
$ws.Range("M122").Value = -19751.7139
$totalCells++
$ws.Range("N122").Value = -46855
$totalCells++

# Row 128
$ws.Range("H128").Value = 41550.625
$totalCells++
$ws.Range("J128").Value = 41550.625
$totalCells++
$ws.Range("L128").Value = 41550.625
$totalCells++
$ws.Range("N128").Value = -51510.625
$totalCells++

# Row 134
$ws.Range("H134").Value = 61626
$totalCells++
$ws.Range("J134").Value = 61626
$totalCells++
$ws.Range("L134").Value = 184878
$totalCells++
$ws.Range("N134").Value = -189948
$totalCells++

Write-Host "Updated $totalCells cells across 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)."
